# Auto-generated Excel COM-interop script to update Asura Profits sheets
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per scheduled market-data refresh
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1166.6666
$ws.Cells.Item(4, 10).Value = 1500
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 14).Value = -1728
$ws.Cells.Item(43, 8).Value = 8747.0625
$ws.Cells.Item(43, 9).Value = 12073.1
$ws.Cells.Item(43, 10).Value = 3203.6667
$ws.Cells.Item(43, 11).Value = 12073.1
$ws.Cells.Item(43, 12).Value = 3203.6667
$ws.Cells.Item(43, 13).Value = -12004.1
$ws.Cells.Item(43, 14).Value = -3341.6667
$ws.Cells.Item(113, 8).Value = 2956.75
$ws.Cells.Item(113, 9).Value = 2338.3333
$ws.Cells.Item(113, 10).Value = 3162.889
$ws.Cells.Item(113, 11).Value = 2338.3333
$ws.Cells.Item(113, 12).Value = 3162.889
$ws.Cells.Item(113, 13).Value = 915.6667000000002
$ws.Cells.Item(113, 14).Value = -9670.888999999999
$ws.Cells.Item(127, 8).Value = 1605.909
$ws.Cells.Item(127, 9).Value = 529.1667
$ws.Cells.Item(127, 11).Value = 1587.5001
$ws.Cells.Item(127, 13).Value = 3372.4999
$ws.Cells.Item(135, 8).Value = 1408.75
$ws.Cells.Item(135, 9).Value = 1032.4
$ws.Cells.Item(135, 10).Value = 2036
$ws.Cells.Item(135, 11).Value = 9291.6
$ws.Cells.Item(135, 12).Value = 18324
$ws.Cells.Item(135, 13).Value = -6756.6
$ws.Cells.Item(135, 14).Value = -23394
$ws.Cells.Item(138, 8).Value = 3706.7344
$ws.Cells.Item(138, 9).Value = 2663.7368
$ws.Cells.Item(138, 10).Value = 4147.1113
$ws.Cells.Item(138, 11).Value = 7991.2104
$ws.Cells.Item(138, 12).Value = 12441.3339
$ws.Cells.Item(138, 13).Value = -2851.2104
$ws.Cells.Item(138, 14).Value = -22721.3339

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(5, 8).Value = 2209.8
$ws.Cells.Item(5, 9).Value = 2209.8
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 2209.8
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -2097.8
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 851.3125
$ws.Cells.Item(22, 9).Value = 676.75
$ws.Cells.Item(22, 11).Value = 676.75
$ws.Cells.Item(22, 13).Value = -377.75
$ws.Cells.Item(32, 8).Value = 12113.328
$ws.Cells.Item(32, 9).Value = 13217.148
$ws.Cells.Item(32, 10).Value = 6152.7
$ws.Cells.Item(32, 11).Value = 13217.148
$ws.Cells.Item(32, 12).Value = 6152.7
$ws.Cells.Item(32, 13).Value = -12930.148
$ws.Cells.Item(32, 14).Value = -6726.7
$ws.Cells.Item(41, 8).Value = 1839.8572
$ws.Cells.Item(41, 9).Value = 1839.8572
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 1839.8572
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -1425.8572
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(88, 8).Value = 2743.7144
$ws.Cells.Item(88, 9).Value = 2235.3333
$ws.Cells.Item(88, 10).Value = 3125
$ws.Cells.Item(88, 11).Value = 2235.3333
$ws.Cells.Item(88, 12).Value = 3125
$ws.Cells.Item(88, 13).Value = -1829.3333
$ws.Cells.Item(88, 14).Value = -3937
$ws.Cells.Item(91, 8).Value = 2743.7144
$ws.Cells.Item(91, 9).Value = 2235.3333
$ws.Cells.Item(91, 10).Value = 3125
$ws.Cells.Item(91, 11).Value = 2235.3333
$ws.Cells.Item(91, 12).Value = 3125
$ws.Cells.Item(91, 13).Value = -831.3332999999998
$ws.Cells.Item(91, 14).Value = -5933
$ws.Cells.Item(122, 8).Value = 5390.478
$ws.Cells.Item(122, 10).Value = 2200
$ws.Cells.Item(122, 12).Value = 6600
$ws.Cells.Item(122, 14).Value = -11500

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 2209.8
$ws.Cells.Item(4, 9).Value = 2209.8
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2209.8
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -2094.8
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 48792.715
$ws.Cells.Item(20, 9).Value = 67655.92999999999
$ws.Cells.Item(20, 10).Value = 1634.6666
$ws.Cells.Item(20, 11).Value = 67655.92999999999
$ws.Cells.Item(20, 12).Value = 1634.6666
$ws.Cells.Item(20, 13).Value = -67408.92999999999
$ws.Cells.Item(20, 14).Value = -2128.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 93.333336
$ws.Cells.Item(7, 9).Value = 93.333336
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 93.333336
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 19.666664
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(127, 8).Value = 25000
$ws.Cells.Item(127, 10).Value = 25000
$ws.Cells.Item(127, 12).Value = 25000
$ws.Cells.Item(127, 14).Value = -34920

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 715.9
$ws.Cells.Item(46, 10).Value = 873.3333
$ws.Cells.Item(46, 12).Value = 2619.9999
$ws.Cells.Item(46, 14).Value = -2801.9999
$ws.Cells.Item(68, 8).Value = 824.5946
$ws.Cells.Item(68, 9).Value = 597.8570999999999
$ws.Cells.Item(68, 10).Value = 1530
$ws.Cells.Item(68, 11).Value = 1793.5713
$ws.Cells.Item(68, 12).Value = 4590
$ws.Cells.Item(68, 13).Value = -982.5712999999998
$ws.Cells.Item(68, 14).Value = -6212
$ws.Cells.Item(71, 8).Value = 824.5946
$ws.Cells.Item(71, 9).Value = 597.8570999999999
$ws.Cells.Item(71, 10).Value = 1530
$ws.Cells.Item(71, 11).Value = 5380.7139
$ws.Cells.Item(71, 12).Value = 13770
$ws.Cells.Item(71, 13).Value = -1324.7139
$ws.Cells.Item(71, 14).Value = -21882
$ws.Cells.Item(107, 8).Value = 838.16
$ws.Cells.Item(107, 9).Value = 396.33334
$ws.Cells.Item(107, 10).Value = 1027.5143
$ws.Cells.Item(107, 11).Value = 1189.00002
$ws.Cells.Item(107, 12).Value = 3082.5429
$ws.Cells.Item(107, 13).Value = 730.9999800000001
$ws.Cells.Item(107, 14).Value = -6922.5429
$ws.Cells.Item(131, 8).Value = 2601.6936
$ws.Cells.Item(131, 10).Value = 3465.5227
$ws.Cells.Item(131, 12).Value = 10396.5681
$ws.Cells.Item(131, 14).Value = -20476.5681
$ws.Cells.Item(138, 8).Value = 2389.0715
$ws.Cells.Item(138, 10).Value = 3850.8572
$ws.Cells.Item(138, 12).Value = 11552.5716
$ws.Cells.Item(138, 14).Value = -21832.5716
$ws.Cells.Item(139, 8).Value = 1706.7931
$ws.Cells.Item(139, 9).Value = 1557.5769
$ws.Cells.Item(139, 10).Value = 3000
$ws.Cells.Item(139, 11).Value = 4672.7307
$ws.Cells.Item(139, 12).Value = 9000
$ws.Cells.Item(139, 13).Value = 467.2692999999999
$ws.Cells.Item(139, 14).Value = -19280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5934.696
$ws.Cells.Item(70, 9).Value = 5899.8945
$ws.Cells.Item(70, 11).Value = 5899.8945
$ws.Cells.Item(70, 13).Value = -5629.8945
$ws.Cells.Item(73, 8).Value = 5934.696
$ws.Cells.Item(73, 9).Value = 5899.8945
$ws.Cells.Item(73, 11).Value = 5899.8945
$ws.Cells.Item(73, 13).Value = -4963.8945
$ws.Cells.Item(80, 8).Value = 2900.2
$ws.Cells.Item(80, 9).Value = 2708.9092
$ws.Cells.Item(80, 10).Value = 3134
$ws.Cells.Item(80, 11).Value = 2708.9092
$ws.Cells.Item(80, 12).Value = 3134
$ws.Cells.Item(80, 13).Value = -1710.9092
$ws.Cells.Item(80, 14).Value = -5130
$ws.Cells.Item(83, 8).Value = 2900.2
$ws.Cells.Item(83, 9).Value = 2708.9092
$ws.Cells.Item(83, 10).Value = 3134
$ws.Cells.Item(83, 11).Value = 13544.546
$ws.Cells.Item(83, 12).Value = 15670
$ws.Cells.Item(83, 13).Value = -8552.546
$ws.Cells.Item(83, 14).Value = -25654
$ws.Cells.Item(102, 8).Value = 3599.8
$ws.Cells.Item(102, 9).Value = 3466.4
$ws.Cells.Item(102, 11).Value = 3466.4
$ws.Cells.Item(102, 13).Value = -1844.4
$ws.Cells.Item(107, 8).Value = 5955.4736
$ws.Cells.Item(107, 9).Value = 7885.7856
$ws.Cells.Item(107, 10).Value = 550.6
$ws.Cells.Item(107, 11).Value = 7885.7856
$ws.Cells.Item(107, 12).Value = 550.6
$ws.Cells.Item(107, 13).Value = -5965.7856
$ws.Cells.Item(107, 14).Value = -4390.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2110.7144
$ws.Cells.Item(68, 9).Value = 1606.25
$ws.Cells.Item(68, 10).Value = 2783.3333
$ws.Cells.Item(68, 11).Value = 1606.25
$ws.Cells.Item(68, 12).Value = 2783.3333
$ws.Cells.Item(68, 13).Value = -857.25
$ws.Cells.Item(68, 14).Value = -4281.3333
$ws.Cells.Item(71, 8).Value = 2110.7144
$ws.Cells.Item(71, 9).Value = 1606.25
$ws.Cells.Item(71, 10).Value = 2783.3333
$ws.Cells.Item(71, 11).Value = 8031.25
$ws.Cells.Item(71, 12).Value = 13916.6665
$ws.Cells.Item(71, 13).Value = -4287.25
$ws.Cells.Item(71, 14).Value = -21404.6665
$ws.Cells.Item(98, 8).Value = 46333.332
$ws.Cells.Item(98, 10).Value = 46333.332
$ws.Cells.Item(98, 12).Value = 46333.332
$ws.Cells.Item(98, 14).Value = -52323.332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1973
$ws.Cells.Item(122, 9).Value = 1850.75
$ws.Cells.Item(122, 10).Value = 2136
$ws.Cells.Item(122, 11).Value = 5552.25
$ws.Cells.Item(122, 12).Value = 6408
$ws.Cells.Item(122, 13).Value = -3102.25
$ws.Cells.Item(122, 14).Value = -11308

Write-Host "Asura Profits sheets updated successfully"